# Update the two-digit / one-digit division answers table.
# Cells are addressed positionally (row, column) because several
# original cell values are duplicated (e.g. "34÷5=6, 4" appears twice
# in row 1) but map to different replacements, so a global
# Find/Replace-All would be incorrect.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "29÷6=4, 5"
$t.Cell(1, 2).Range.Text = "30÷8=3, 6"
$t.Cell(1, 3).Range.Text = "82÷5=16, 2"
$t.Cell(1, 4).Range.Text = "63÷9=7, 0"
$t.Cell(1, 5).Range.Text = "92÷6=15, 2"
$t.Cell(5, 1).Range.Text = "25÷9=2, 7"
$t.Cell(5, 2).Range.Text = "56÷6=9, 2"
$t.Cell(5, 3).Range.Text = "29÷7=4, 1"
$t.Cell(5, 4).Range.Text = "92÷8=11, 4"
$t.Cell(5, 5).Range.Text = "57÷9=6, 3"
$t.Cell(9, 1).Range.Text = "31÷4=7, 3"
$t.Cell(9, 2).Range.Text = "24÷6=4, 0"
$t.Cell(9, 3).Range.Text = "59÷6=9, 5"
$t.Cell(9, 4).Range.Text = "21÷7=3, 0"
$t.Cell(9, 5).Range.Text = "65÷6=10, 5"
$t.Cell(13, 1).Range.Text = "31÷7=4, 3"
$t.Cell(13, 2).Range.Text = "22÷6=3, 4"
$t.Cell(13, 3).Range.Text = "79÷2=39, 1"
$t.Cell(13, 4).Range.Text = "62÷6=10, 2"
$t.Cell(13, 5).Range.Text = "81÷4=20, 1"
$t.Cell(17, 1).Range.Text = "20÷4=5, 0"
$t.Cell(17, 2).Range.Text = "85÷7=12, 1"
$t.Cell(17, 3).Range.Text = "27÷7=3, 6"
$t.Cell(17, 4).Range.Text = "56÷8=7, 0"
$t.Cell(17, 5).Range.Text = "10÷3=3, 1"

Write-Host "Done updating division answers table."
